$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Database fix: Javelin (row 22) range "1~1" -> "1~2" ---
$ws.Range("D22").Value = 2

# --- Column B width, set so saved cars/units render with their full name ---
$ws.Columns.Item(2).ColumnWidth = 17.33

# --- Scroll/selection: no longer pinned near the bottom of the sheet, and
#     the active cell/selection moved from L2 to E17 ---
$ws.Range("E17").Select() | Out-Null

Write-Output "Edits applied"
